$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# "Almacenamiento_Producto" table (rows 7-10) gains a new attribute column
# "totalUnidadesOfrecidas", inserted as the new column I. Everything that used
# to live in columns J..M (the table's id/idContenedorActual/
# existenciasActuales/nivelReorden columns) shifts one column right, to K..N.
# ---------------------------------------------------------------------------

function Copy-Cell($srcRow, $srcCol, $dstRow, $dstCol) {
    $src = $ws.Cells.Item($srcRow, $srcCol)
    $dst = $ws.Cells.Item($dstRow, $dstCol)
    $src.Copy($dst)
}

# Row 7 only has the table title in column J -> moves to column K.
Copy-Cell 7 10 7 11
$ws.Cells.Item(7, 10).Clear()

# Rows 8 (headers), 9 (values), 10 (blank separator row) each have J..M
# populated; shift them right into K..N (right-to-left so we never
# overwrite a cell we still need to read).
foreach ($r in 8..10) {
    Copy-Cell $r 13 $r 14   # M -> N
    Copy-Cell $r 12 $r 13   # L -> M
    Copy-Cell $r 11 $r 12   # K -> L
    Copy-Cell $r 10 $r 11   # J -> K
    $ws.Cells.Item($r, 10).Clear()   # clear the now-stale column J cell
}

# New column I content for the "Almacenamiento_Producto" table: header,
# value placeholder and blank separator, matching the formatting already
# used by the neighbouring "idProducto" column (H).
$ws.Cells.Item(8, 8).Copy($ws.Cells.Item(8, 9))
$ws.Cells.Item(8, 9).Value = "totalUnidadesOfrecidas"

$ws.Cells.Item(9, 8).Copy($ws.Cells.Item(9, 9))
$ws.Cells.Item(9, 9).Value = "DD"

$ws.Cells.Item(10, 8).Copy($ws.Cells.Item(10, 9))
$ws.Cells.Item(10, 9).ClearContents()

# ---------------------------------------------------------------------------
# Unrelated cosmetic row-height tweaks that came along with the same edit.
# ---------------------------------------------------------------------------
$ws.Rows.Item(34).RowHeight = 54
$ws.Rows.Item(44).RowHeight = 36

# ---------------------------------------------------------------------------
# View state: zoom level and current selection.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 63
$ws.Range("J9").Select() | Out-Null
